# Insert a new data row above the current row 33 ("Haba" price list for
# Terminal Hortofrutícola Agro Chillán). This pushes the existing rows
# 33:67 down to 34:68, preserving every cell (values + date style) exactly
# as they were - which matches the diff's row-by-row shift pattern.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("33:33").Insert()

# Populate the newly-inserted row 33 with its own (new) record.
$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(33, 3).Value = "Ñuble"
$ws.Cells.Item(33, 4).Value = 44895
$ws.Cells.Item(33, 5).Value = 16
$ws.Cells.Item(33, 6).Value = 100112026
$ws.Cells.Item(33, 7).Value = "Haba"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 60
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 12000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(33, 16).Value = 480
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
